$wb = $excel.ActiveWorkbook

# --- Update the descriptive text on "Hoja1" (sheet 1), cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.41 = 30368.59 pesos`n✅ 30368.59 pesos = 7.39 = 934.75 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update numeric rate values on "tasas" (sheet 2) ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 134.975
$ws2.Range("O10").Value = 4099
$ws2.Range("N12").Value = 4109.8
$ws2.Range("O12").Value = 126.5
